$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price-column updates whose new text would otherwise be auto-parsed as a
# number by Excel (losing the trailing zero / exact digit count, e.g. "0.600"
# -> 0.6). Force the cell to Text format first so the literal string from the
# source data lands in the sheet unchanged, exactly like the original
# inline-string cell.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

# Apply the updated Price (D) and Volume(1h) (E) values for each changed row.
$ws.Range("D2").Value = "63.351.85"
$ws.Range("E2").Value = "  +1.98%  "
$ws.Range("D3").Value = "2.471.80"
$ws.Range("E3").Value = "  +2.66%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "565.38"
$ws.Range("E5").Value = "  +0.81%  "
$ws.Range("D6").Value = "143.41"
$ws.Range("E6").Value = "  +3.58%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "0.586"
$ws.Range("E8").Value = "  -0.39%  "
$ws.Range("D9").Value = "2.472.04"
$ws.Range("E9").Value = "  +2.75%  "
$ws.Range("E10").Value = "  +0.66%  "
$ws.Range("D11").Value = "5.74"
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("E12").Value = "  +1.58%  "
$ws.Range("E13").Value = "  +2.06%  "
$ws.Range("D14").Value = "27.36"
$ws.Range("E14").Value = "  +5.98%  "
$ws.Range("D15").Value = "2.910.05"
$ws.Range("E15").Value = "  +2.75%  "
$ws.Range("D16").Value = "63.117.08"
$ws.Range("E16").Value = "  +1.72%  "
$ws.Range("E17").Value = "  +2.52%  "
$ws.Range("D18").Value = "2.465.61"
$ws.Range("E18").Value = "  +2.28%  "
$ws.Range("D19").Value = "11.30"
$ws.Range("E19").Value = "  +2.43%  "
$ws.Range("D20").Value = "341.30"
$ws.Range("E20").Value = "  -0.75%  "
$ws.Range("E21").Value = "  +1.61%  "
$ws.Range("E22").Value = "  -1.66%  "
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").Value = "65.62"
$ws.Range("E24").Value = "  +0.72%  "
$ws.Range("E25").Value = "  -1.12%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  +0.29%  "
$ws.Range("E28").Value = "  +4.74%  "
$ws.Range("E29").Value = "  -3.24%  "
$ws.Range("D30").Value = "1.86"
$ws.Range("E30").Value = "  +2.46%  "
$ws.Range("E31").Value = "  +5.85%  "
$ws.Range("D32").Value = "0.0₃0798"
$ws.Range("E32").Value = "  +3.06%  "
$ws.Range("D33").Value = "176.37"
$ws.Range("E33").Value = "  +2.72%  "
$ws.Range("E34").Value = "  +7.54%  "
$ws.Range("D35").Value = "398.71"
$ws.Range("E35").Value = "  +9.82%  "
$ws.Range("E36").Value = "  +1.18%  "
$ws.Range("D37").Value = "18.86"
$ws.Range("E37").Value = "  +1.69%  "
$ws.Range("E39").Value = "  -4.06%  "
$ws.Range("D40").Value = "1.75"
$ws.Range("E40").Value = "  +4.50%  "
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").Value = "40.47"
$ws.Range("E42").Value = "  +3.86%  "
$ws.Range("D43").Value = "150.24"
$ws.Range("E43").Value = "  +4.26%  "
$ws.Range("D44").Value = "3.72"
$ws.Range("E44").Value = "  +1.14%  "
$ws.Range("D45").Value = "20.68"
$ws.Range("E45").Value = "  +0.25%  "
$ws.Range("D46").Value = "0.600"
$ws.Range("E46").Value = "  +2.90%  "
$ws.Range("E47").Value = "  -0.12%  "
$ws.Range("D48").Value = "0.0518"
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("D49").Value = "0.0229"
$ws.Range("E49").Value = "  +3.25%  "
$ws.Range("E50").Value = "  +1.00%  "
$ws.Range("E51").Value = "  +5.54%  "
